$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain text while writing (otherwise Excel
# auto-converts numeric-looking strings to floating point numbers and
# mangles values like "1.773.38" or loses trailing zeros like "0.0600").
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "25.212.88"
$ws.Range("E2").Value = "  -2.60%  "

$ws.Range("D3").Value = "1.552.79"
$ws.Range("E3").Value = "  -4.18%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "205.78"
$ws.Range("E5").Value = "  -3.59%  "

$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").Value = "0.473"
$ws.Range("E7").Value = "  -5.56%  "

$ws.Range("D8").Value = "0.0600"
$ws.Range("E8").Value = "  -2.20%  "

$ws.Range("E9").Value = "  -3.65%  "

$ws.Range("D10").Value = "17.62"
$ws.Range("E10").Value = "  -2.99%  "

$ws.Range("D11").Value = "0.0778"
$ws.Range("E11").Value = "  -0.96%  "

$ws.Range("D12").Value = "1.773.38"
$ws.Range("E12").Value = "  -3.97%  "

$ws.Range("D13").Value = "1.553.58"
$ws.Range("E13").Value = "  -4.19%  "

$ws.Range("D14").Value = "3.94"
$ws.Range("E14").Value = "  -5.30%  "

$ws.Range("D15").Value = "0.498"
$ws.Range("E15").Value = "  -4.69%  "

$ws.Range("D16").Value = "25.192.50"
$ws.Range("E16").Value = "  -2.70%  "

$ws.Range("D17").Value = "0.0₃0703"
$ws.Range("E17").Value = "  -3.70%  "

$ws.Range("D18").Value = "58.38"
$ws.Range("E18").Value = "  -4.26%  "

$ws.Range("D19").Value = "1.01"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").Value = "183.46"
$ws.Range("E20").Value = "  -4.36%  "

$ws.Range("D21").Value = "4.07"
$ws.Range("E21").Value = "  -3.48%  "

$ws.Range("D22").Value = "9.16"
$ws.Range("E22").Value = "  -3.83%  "

$ws.Range("D23").Value = "5.80"
$ws.Range("E23").Value = "  -4.11%  "

$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").Value = "0.126"
$ws.Range("E25").Value = "  -4.30%  "

$ws.Range("D26").Value = "139.13"
$ws.Range("E26").Value = "  -3.13%  "

$ws.Range("E27").Value = "  -4.90%  "

$ws.Range("D28").Value = "14.69"
$ws.Range("E28").Value = "  -2.75%  "

$ws.Range("D29").Value = "6.36"
$ws.Range("E29").Value = "  -5.04%  "

$ws.Range("D30").Value = "1.14"
$ws.Range("E30").Value = "  -6.72%  "

$ws.Range("D31").Value = "0.0460"
$ws.Range("E31").Value = "  -4.33%  "

$ws.Range("D32").Value = "2.98"
$ws.Range("E32").Value = "  -3.74%  "

$ws.Range("D33").Value = "2.94"
$ws.Range("E33").Value = "  -4.89%  "

$ws.Range("D34").Value = "1.43"
$ws.Range("E34").Value = "  -3.48%  "

$ws.Range("E35").Value = "  -4.02%  "

$ws.Range("D36").Value = "1.079.73"
$ws.Range("E36").Value = "  -3.14%  "

$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("D38").Value = "0.0148"
$ws.Range("E38").Value = "  -2.73%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.26"
$ws.Range("E39").Value = "  -7.12%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.808"
$ws.Range("E40").Value = "  +5.64%  "

$ws.Range("D41").Value = "0.488"
$ws.Range("E41").Value = "  -5.40%  "

$ws.Range("D42").Value = "0.750"
$ws.Range("E42").Value = "  -10.92%  "

$ws.Range("D43").Value = "92.16"
$ws.Range("E43").Value = "  -5.82%  "

$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.686.78"
$ws.Range("E44").Value = "  -4.03%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "4.99"
$ws.Range("E45").Value = "  -3.23%  "

$ws.Range("D46").Value = "0.0₆0107"
$ws.Range("E46").Value = "  -7.08%  "

$ws.Range("D47").Value = "51.93"
$ws.Range("E47").Value = "  -4.14%  "

$ws.Range("D48").Value = "0.0503"
$ws.Range("E48").Value = "  -5.00%  "

$ws.Range("D49").Value = "1.42"
$ws.Range("E49").Value = "  -2.31%  "

$ws.Range("D50").Value = "0.406"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("E51").Value = "  -0.23%  "

# Restore the default cell style on the Price column so the written cells
# keep no explicit style index, matching the original file formatting.
$dRange.Style = "Normal"
